$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell E1 ("Billing inv") -----------------------------------
# Same look as A1 (Meeting): bold font, green fill (FF00B050), medium border,
# centered alignment.
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Billing inv"

# --- New header cell F1 ("Receivables") -----------------------------------
# Same look as D1 (Billing req: bold, bordered, centered) but with a new
# light-blue fill color (FF00B0F0).
$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Interior.Color = 15773696
$ws.Range("F1").Value = "Receivables"

# --- Updated data row 2 -----------------------------------------------------
$ws.Range("A2").Value = "ME-634"
$ws.Range("B2").Value = "OF-639"
$ws.Range("C2").Value = "CT-305"
$ws.Range("D2").Value = "BR-419"

$ws.Range("E2").Value = "BI-143"
$ws.Range("E2").Style = "Normal"

$ws.Range("F2").Value = "RE-69"

# --- Column widths for the two new columns ----------------------------------
$ws.Columns.Item(5).ColumnWidth = 13.6667
$ws.Columns.Item(6).ColumnWidth = 10.6667

# --- Selection moves to E9 --------------------------------------------------
$ws.Range("E9").Select()
